$d = $word.ActiveDocument

# The "KEY ACHIEVEMENTS AND IMPACT" section contains several bullets whose wording is
# identical (or near-identical) to bullets that already appear earlier, under
# "PROFESSIONAL EXPERIENCE" (Siege Analytics). To avoid touching those earlier
# occurrences, every lookup below re-derives the live start/end of the achievements
# section (paragraph offsets shift as we edit) and scopes the work to it.

function Get-AchievementsSectionRange() {
    $start = -1
    $end = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($start -eq -1) {
            if ($t -match "KEY ACHIEVEMENTS AND IMPACT") {
                $start = $p.Range.Start
            }
            continue
        }
        if ($t -match "TECHNICAL SKILLS") {
            $end = $p.Range.Start
            break
        }
    }
    if ($start -eq -1) { $start = 0 }
    if ($end -eq -1) { $end = $d.Content.End }
    return $d.Range($start, $end)
}

function Replace-InAchievements([string]$oldText, [string]$newText) {
    $r = Get-AchievementsSectionRange
    $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

function Remove-ParagraphInAchievements([string]$containsText) {
    $r = Get-AchievementsSectionRange
    $start = $r.Start
    $end = $r.End
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -ge $start -and $p.Range.Start -lt $end -and $p.Range.Text -match [regex]::Escape($containsText)) {
            $p.Range.Delete()
            return
        }
    }
}

# 1) "Achieved 87% prediction accuracy..." -> "Revenue generation: Delivered $4.9M additional revenue through optimization"
Replace-InAchievements `
    "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%" `
    "Revenue generation: Delivered `$4.9M additional revenue through optimization"

# 2) "Delivered $4.9M additional revenue through continuous testing and optimization, increased conversion rates by 23%" -> "23% conversion rate improvement"
Replace-InAchievements `
    "Delivered `$4.9M additional revenue through continuous testing and optimization, increased conversion rates by 23%" `
    "23% conversion rate improvement"

# 3) "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations"
#    -> "Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
Replace-InAchievements `
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations" `
    "Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"

# 4) Remove the whole "Developed longitudinal data analysis methods..." bullet paragraph (incl. its mark)
Remove-ParagraphInAchievements "Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality"

# 5) Remove the whole "Discovered systematic race coding errors..." bullet paragraph
Remove-ParagraphInAchievements "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%"

# 6) "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis"
#    -> "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
Replace-InAchievements `
    "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis" `
    "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

Write-Output "Done."
